$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 87 (shifts old row 87 -> 88, 88 -> 89, etc.)
$ws.Rows.Item(87).Insert()

# Row 86: change end time value
$ws.Range("E86").Value = 0.875

# Row 87: fill in new data row
$ws.Range("A87").Value = 2014
$ws.Range("B87").Value = 3
$ws.Range("C87").Value = 21
$ws.Range("D87").Value = 0.88888888888888884
$ws.Range("E87").Value = 0.91666666666666663
$ws.Range("F87").Formula = "=(E87-D87)*24*60"
$ws.Range("G87").Formula = "=F87/60"

# Update selection to match new target
$ws.Range("A88").Select()
